$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from AC1 to the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill data rows 2-58 with team record values
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 30).Value = 66   # AD
    $ws.Cells.Item($r, 31).Value = 96   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
